# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff:
# "Updated cryptos list on Thu Apr 11 07:56:29 UTC 2024 with GitHub Actions"
#
# Cells in column D that contain plain numeric-looking text (e.g. "605.06")
# must be forced to Text format first, otherwise Excel's COM layer silently
# coerces the assigned string into a floating point number (losing the
# original text representation / exact digits), same as would happen if you
# typed the digits into a General-formatted cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.909.82'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '3.587.80'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.06'
$ws.Range('E5').Value = '  +3.03%  '
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('D7').Value = '3.582.70'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.201'
$ws.Range('E10').Value = '  +5.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.48'
$ws.Range('E11').Value = '  +10.18%  '
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.34'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').Value = '4.169.48'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '621.82'
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('D18').Value = '3.589.95'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('D19').Value = '71.079.19'
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('E20').Value = '  -2.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.54'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.38'
$ws.Range('E23').Value = '  -16.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.25'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.17'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.67'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.39'
$ws.Range('E29').Value = '  +4.52%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.38'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.09'
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.22'
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.32'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '629.23'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.83'
$ws.Range('E36').Value = '  +9.19%  '
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.91'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('E39').Value = '  +6.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '57.54'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.143'
$ws.Range('E42').Value = '  +4.21%  '
$ws.Range('D43').Value = '3.419.60'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.327'
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.04'
$ws.Range('E45').Value = '  +10.48%  '
$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = '0.0₃0720'
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.72'
$ws.Range('E47').Value = '  +5.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.18'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.82'
$ws.Range('E50').Value = '  +0.22%  '
